{"js": "// Update the worksheet date line and every two-digit x two-digit\n// multiplication problem/answer in the table to the new set of values.\n// We match each old text exactly (matchCase, not a wildcard) and replace\n// the whole run's text with the new value, preserving all formatting\n// because we operate on the existing run's range rather than inserting a\n// brand new paragraph/run.\n\nconst replacements = [\n  [\"2024-08-16 Friday\", \"2024-08-17 Saturday\"],\n  [\"51\u00d777=3927\", \"79\u00d720=1580\"],\n  [\"38\u00d792=3496\", \"24\u00d733=792\"],\n  [\"83\u00d748=3984\", \"61\u00d724=1464\"],\n  [\"26\u00d776=1976\", \"89\u00d712=1068\"],\n  [\"85\u00d771=6035\", \"58\u00d756=3248\"],\n  [\"36\u00d734=1224\", \"33\u00d798=3234\"],\n  [\"89\u00d719=1691\", \"36\u00d743=1548\"],\n  [\"72\u00d749=3528\", \"92\u00d782=7544\"],\n  [\"26\u00d719=494\", \"77\u00d789=6853\"],\n  [\"61\u00d799=6039\", \"93\u00d799=9207\"],\n  [\"74\u00d727=1998\", \"76\u00d711=836\"],\n  [\"98\u00d732=3136\", \"42\u00d774=3108\"],\n  [\"66\u00d754=3564\", \"84\u00d739=3276\"],\n  [\"74\u00d787=6438\", \"32\u00d754=1728\"],\n  [\"93\u00d752=4836\", \"40\u00d726=1040\"],\n  [\"72\u00d779=5688\", \"84\u00d735=2940\"],\n  [\"68\u00d716=1088\", \"72\u00d742=3024\"],\n  [\"80\u00d711=880\", \"43\u00d741=1763\"],\n  [\"67\u00d746=3082\", \"84\u00d757=4788\"],\n  [\"59\u00d771=4189\", \"65\u00d745=2925\"],\n  [\"23\u00d724=552\", \"80\u00d754=4320\"],\n  [\"70\u00d751=3570\", \"64\u00d751=3264\"],\n  [\"91\u00d715=1365\", \"43\u00d733=1419\"],\n  [\"97\u00d756=5432\", \"31\u00d778=2418\"],\n  [\"38\u00d714=532\", \"78\u00d726=2028\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date line and every two-digit x two-digit\n# multiplication problem/answer in the table to the new set of values.\n# Each old->new pair is applied with Find/Replace (wdReplaceAll semantics)\n# against the whole document body so formatting on the existing runs is\n# left untouched; only the literal text is rewritten.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-08-16 Friday\", \"2024-08-17 Saturday\"),\n  @(\"51\u00d777=3927\", \"79\u00d720=1580\"),\n  @(\"38\u00d792=3496\", \"24\u00d733=792\"),\n  @(\"83\u00d748=3984\", \"61\u00d724=1464\"),\n  @(\"26\u00d776=1976\", \"89\u00d712=1068\"),\n  @(\"85\u00d771=6035\", \"58\u00d756=3248\"),\n  @(\"36\u00d734=1224\", \"33\u00d798=3234\"),\n  @(\"89\u00d719=1691\", \"36\u00d743=1548\"),\n  @(\"72\u00d749=3528\", \"92\u00d782=7544\"),\n  @(\"26\u00d719=494\", \"77\u00d789=6853\"),\n  @(\"61\u00d799=6039\", \"93\u00d799=9207\"),\n  @(\"74\u00d727=1998\", \"76\u00d711=836\"),\n  @(\"98\u00d732=3136\", \"42\u00d774=3108\"),\n  @(\"66\u00d754=3564\", \"84\u00d739=3276\"),\n  @(\"74\u00d787=6438\", \"32\u00d754=1728\"),\n  @(\"93\u00d752=4836\", \"40\u00d726=1040\"),\n  @(\"72\u00d779=5688\", \"84\u00d735=2940\"),\n  @(\"68\u00d716=1088\", \"72\u00d742=3024\"),\n  @(\"80\u00d711=880\", \"43\u00d741=1763\"),\n  @(\"67\u00d746=3082\", \"84\u00d757=4788\"),\n  @(\"59\u00d771=4189\", \"65\u00d745=2925\"),\n  @(\"23\u00d724=552\", \"80\u00d754=4320\"),\n  @(\"70\u00d751=3570\", \"64\u00d751=3264\"),\n  @(\"91\u00d715=1365\", \"43\u00d733=1419\"),\n  @(\"97\u00d756=5432\", \"31\u00d778=2418\"),\n  @(\"38\u00d714=532\", \"78\u00d726=2028\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
